$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# A leading apostrophe forces Excel to treat numeric-looking values (e.g. "588.97")
# as text, matching the original inlineStr cell type; Style is reset to "Normal"
# afterwards so no stray text-format style gets attached to the cell.

$ws.Range("D2").Value = "'63.026.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.64%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.144.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'588.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.63%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.73%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.18%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.139.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.91%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.97%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -2.98%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'34.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.07%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.663.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.13%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.87%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.146.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.10%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'63.011.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.87%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'474.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.28%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.39%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.84%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.33%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'84.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.82%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'12.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.77%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.01%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.29%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.37%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.99%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.05%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'26.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.31%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -4.86%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.37%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.83%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'52.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.35%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0696"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -7.11%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.10%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'414.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -6.15%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.920.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -6.61%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.03%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -4.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'25.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.67%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -8.00%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'120.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.71%  "
$ws.Range("E51").Style = "Normal"
